$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40: Stuck in the Moment | Horn Glue
$ws.Range("H40").Value = 1804.15
$ws.Range("I40").Value = 1844.7778
$ws.Range("J40").Value = 1770.909
$ws.Range("K40").Value = 1844.7778
$ws.Range("L40").Value = 1770.909
$ws.Range("M40").Value = -1669.7778
$ws.Range("N40").Value = -2120.909

# Row 74: Adhesive of Antipathy | Wing Glue
$ws.Range("H74").Value = 3948.2
$ws.Range("I74").Value = 3873
$ws.Range("J74").Value = 3998.3333
$ws.Range("K74").Value = 3873
$ws.Range("L74").Value = 3998.3333
$ws.Range("M74").Value = -2937
$ws.Range("N74").Value = -5870.3333

# Row 77: It's Gonna Grow Back (L) | Wing Glue
$ws.Range("H77").Value = 3948.2
$ws.Range("I77").Value = 3873
$ws.Range("J77").Value = 3998.3333
$ws.Range("K77").Value = 19365
$ws.Range("L77").Value = 19991.6665
$ws.Range("M77").Value = -14685
$ws.Range("N77").Value = -29351.6665

# Row 100: Asking for a Friend | Beetle Glue
$ws.Range("H100").Value = 2417.611
$ws.Range("I100").Value = 1791.7
$ws.Range("K100").Value = 1791.7
$ws.Range("M100").Value = -1250.7

# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 4399.516
$ws.Range("I132").Value = 4379.1304
$ws.Range("K132").Value = 13137.3912
$ws.Range("M132").Value = -10607.3912

# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 4652294
$ws.Range("I137").Value = 762.23334
$ws.Range("J137").Value = 15386597
$ws.Range("K137").Value = 2286.70002
$ws.Range("L137").Value = 46159791
$ws.Range("M137").Value = 263.2999799999998
$ws.Range("N137").Value = -46164891

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 14289197
$ws.Range("I61").Value = 18521674
$ws.Range("J61").Value = 4588
$ws.Range("K61").Value = 18521674
$ws.Range("L61").Value = 4588
$ws.Range("M61").Value = -18521462
$ws.Range("N61").Value = -5012

# Row 88: The Mast Chance | Adamantite Rivets
$ws.Range("H88").Value = 2470.2354
$ws.Range("I88").Value = 2259
$ws.Range("J88").Value = 2772
$ws.Range("K88").Value = 2259
$ws.Range("L88").Value = 2772
$ws.Range("M88").Value = -1853
$ws.Range("N88").Value = -3584

# Row 91: The Rose and the Riveter (L) | Adamantite Rivets
$ws.Range("H91").Value = 2470.2354
$ws.Range("I91").Value = 2259
$ws.Range("J91").Value = 2772
$ws.Range("K91").Value = 2259
$ws.Range("L91").Value = 2772
$ws.Range("M91").Value = -855
$ws.Range("N91").Value = -5580

# Row 102: Smells of Rich Tama-hagane | Tama-hagane Ingot
$ws.Range("H102").Value = 4000
$ws.Range("J102").Value = 4000
$ws.Range("L102").Value = 4000
$ws.Range("N102").Value = -7244

# Row 122: Haste for High Durium | High Durium Nugget
$ws.Range("H122").Value = 43674.668
$ws.Range("I122").Value = 43674.668
$ws.Range("K122").Value = 131024.004
$ws.Range("M122").Value = -128574.004

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 14289197
$ws.Range("I136").Value = 18521674
$ws.Range("J136").Value = 4588
$ws.Range("K136").Value = 55565022
$ws.Range("L136").Value = 13764
$ws.Range("M136").Value = -55562472
$ws.Range("N136").Value = -18864

$ws = $wb.Worksheets.Item("BSM")
# Row 22: Riveting Run | Iron Rivets
$ws.Range("H22").Value = 395.875
$ws.Range("I22").Value = 386
$ws.Range("K22").Value = 386
$ws.Range("M22").Value = -213

# Row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 22729910
$ws.Range("I86").Value = 3007.1428
$ws.Range("K86").Value = 3007.1428
$ws.Range("M86").Value = -1884.1428

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 22729910
$ws.Range("I89").Value = 3007.1428
$ws.Range("K89").Value = 15035.714
$ws.Range("M89").Value = -9419.714

# Row 103: The Bigger the Blade | Doman Steel Tachi
$ws.Range("H103").Value = 25317
$ws.Range("J103").Value = 25317
$ws.Range("L103").Value = 25317
$ws.Range("N103").Value = -27661

# Row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Range("H105").Value = 3969
$ws.Range("I105").Value = 2597.7778
$ws.Range("J105").Value = 5090.909
$ws.Range("K105").Value = 2597.7778
$ws.Range("L105").Value = 5090.909
$ws.Range("M105").Value = -850.7777999999998
$ws.Range("N105").Value = -8584.909

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 3087.484
$ws.Range("I134").Value = 2220.7273
$ws.Range("J134").Value = 5206.222
$ws.Range("K134").Value = 6662.1819
$ws.Range("L134").Value = 15618.666
$ws.Range("M134").Value = -4127.1819
$ws.Range("N134").Value = -20688.666

$ws = $wb.Worksheets.Item("CRP")
# Row 62: Splinter in the Sewers | Cedar Lumber
$ws.Range("H62").Value = 2298.4614
$ws.Range("I62").Value = 2298.4614
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2298.4614
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -1674.4614

# Row 65: The Lumber of Their Discontent (L) | Cedar Lumber
$ws.Range("H65").Value = 2298.4614
$ws.Range("I65").Value = 2298.4614
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 11492.307
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -8372.307000000001

# Row 93: Reeling for Rods | Muudhorn Fishing Rod
$ws.Range("H93").Value = 12357.625
$ws.Range("I93").Value = 9630.571
$ws.Range("J93").Value = 31447
$ws.Range("K93").Value = 9630.571
$ws.Range("L93").Value = 31447
$ws.Range("M93").Value = -7758.571
$ws.Range("N93").Value = -35191

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 1805.4318
$ws.Range("I132").Value = 1220.7878
$ws.Range("J132").Value = 3559.3635
$ws.Range("K132").Value = 3662.3634
$ws.Range("L132").Value = 10678.0905
$ws.Range("M132").Value = -1132.3634
$ws.Range("N132").Value = -15738.0905

$ws = $wb.Worksheets.Item("CUL")
# Row 3: Trout Fishing in Limsa | Grilled Trout
$ws.Range("H3").Value = 4223.7393
$ws.Range("I3").Value = 2875.875
$ws.Range("J3").Value = 7304.5713
$ws.Range("K3").Value = 8627.625
$ws.Range("L3").Value = 21913.7139
$ws.Range("M3").Value = -8515.625
$ws.Range("N3").Value = -22137.7139

# Row 15: Pretty Enough to Eat | Grilled Carp
$ws.Range("H15").Value = 1805.7142
$ws.Range("I15").Value = 580
$ws.Range("K15").Value = 1740
$ws.Range("M15").Value = -1600

# Row 96: Hunger Is No Game | Popoto Soba
$ws.Range("H96").Value = 5792
$ws.Range("J96").Value = 6722.6665
$ws.Range("L96").Value = 20167.9995
$ws.Range("N96").Value = -24285.9995

# Row 118: Teetotally | Masala Chai
$ws.Range("H118").Value = 1413.8334
$ws.Range("J118").Value = 990
$ws.Range("L118").Value = 2970
$ws.Range("N118").Value = -5456

# Row 134: Don't Knock It Till You've Tried It | Mezcal-marinated Swampmonk
$ws.Range("H134").Value = 5681.1763
$ws.Range("I134").Value = 4330
$ws.Range("J134").Value = 6418.1816
$ws.Range("K134").Value = 12990
$ws.Range("L134").Value = 19254.5448
$ws.Range("M134").Value = -7920
$ws.Range("N134").Value = -29394.5448

# Row 141: Ocean Explosion | Acqua Pazza
$ws.Range("H141").Value = 5149.1177
$ws.Range("I141").Value = 5049
$ws.Range("J141").Value = 5900
$ws.Range("K141").Value = 15147
$ws.Range("L141").Value = 17700
$ws.Range("M141").Value = -9967
$ws.Range("N141").Value = -28060

$ws = $wb.Worksheets.Item("LTW")
# Row 100: Tiger in the Sack | Tiger Leather
$ws.Range("H100").Value = 1760.1666
$ws.Range("I100").Value = 1672.2
$ws.Range("K100").Value = 1672.2
$ws.Range("M100").Value = -1131.2

# Row 122: Hell on Leather | Gaja Leather
$ws.Range("H122").Value = 7067.5
$ws.Range("I122").Value = 9560
$ws.Range("J122").Value = 5643.2144
$ws.Range("K122").Value = 28680
$ws.Range("L122").Value = 16929.6432
$ws.Range("M122").Value = -26230
$ws.Range("N122").Value = -21829.6432

# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 10212053
$ws.Range("I132").Value = 5925.4814
$ws.Range("J132").Value = 22737756
$ws.Range("K132").Value = 17776.4442
$ws.Range("L132").Value = 68213268
$ws.Range("M132").Value = -15246.4442
$ws.Range("N132").Value = -68218328

$ws = $wb.Worksheets.Item("WVR")
# Row 62: Pride Up in Smoke | Rainbow Cloth
$ws.Range("H62").Value = 15786.143
$ws.Range("I62").Value = 6100
$ws.Range("J62").Value = 40001.5
$ws.Range("K62").Value = 6100
$ws.Range("L62").Value = 40001.5
$ws.Range("M62").Value = -5476
$ws.Range("N62").Value = -41249.5

# Row 65: Desperate for Diversionaries (L) | Rainbow Cloth
$ws.Range("H65").Value = 15786.143
$ws.Range("I65").Value = 6100
$ws.Range("J65").Value = 40001.5
$ws.Range("K65").Value = 30500
$ws.Range("L65").Value = 200007.5
$ws.Range("M65").Value = -27380
$ws.Range("N65").Value = -206247.5

# Row 96: Skills on Display | Ruby Cotton Cloth
$ws.Range("H96").Value = 2546
$ws.Range("I96").Value = 1630.5454
$ws.Range("J96").Value = 3138.353
$ws.Range("K96").Value = 1630.5454
$ws.Range("L96").Value = 3138.353
$ws.Range("M96").Value = -257.5454
$ws.Range("N96").Value = -5884.353

# Row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 2285.8
$ws.Range("I126").Value = 1663
$ws.Range("J126").Value = 4777
$ws.Range("K126").Value = 4989
$ws.Range("L126").Value = 14331
$ws.Range("M126").Value = -2519
$ws.Range("N126").Value = -19271

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 1861.95
$ws.Range("I132").Value = 900.8182
$ws.Range("J132").Value = 3036.6667
$ws.Range("K132").Value = 2702.4546
$ws.Range("L132").Value = 9110.000100000001
$ws.Range("M132").Value = -172.4546
$ws.Range("N132").Value = -14170.0001
